$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.127
$ws.Range("E2").Value = 0.014395
$ws.Range("G2").Value = 0.1266191929467616
$ws.Range("H2").Value = 0.1266191929467616
$ws.Range("I2").Value = 0.1244972872160054
$ws.Range("J2").Value = 0.1196481567873167
$ws.Range("K2").Value = 124.83
$ws.Range("L2").Value = 0.105824008138352
$ws.Range("M2").Value = 59.03
$ws.Range("N2").Value = 0.05188264660385319
$ws.Range("O2").Value = 0.4728831210446207
$ws.Range("P2").Value = 59.03
$ws.Range("Q2").Value = 0.05188264660385319
$ws.Range("R2").Value = 0.4728831210446207
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 404.049
$ws.Range("V2").Value = 0.3551267402615666
$ws.Range("W2").Value = 0.1072152999130687
$ws.Range("X2").Value = 0.04500067940374067
$ws.Range("Y2").Value = 0.06221462050932799
$ws.Range("Z2").Value = 1.167080725062257
$ws.Range("AA2").Value = 0.1110655771493076
$ws.Range("AB2").Value = 0.04500067940374067
$ws.Range("AC2").Value = 0.06606489774556697
$ws.Range("AD2").Value = 146.2
$ws.Range("AF2").Value = 146.2
$ws.Range("AG2").Value = -257.849
$ws.Range("AH2").Value = 0.1138664755911399
$ws.Range("AI2").Value = 0.1001507055761063
$ws.Range("AJ2").Value = -0.2930398642589989
$ws.Range("AK2").Value = -0.2442327783729307
$ws.Range("AL2").Value = 9.9
$ws.Range("AM2").Value = 9.9
$ws.Range("AN2").Value = 0.9464620961999093
$ws.Range("AO2").Value = 14.8340404040404
$ws.Range("AP2").Value = -1.669249692496925
$ws.Range("AQ2").Value = 14.8340404040404

# Row 3
$ws.Range("D3").Value = 0.127
$ws.Range("E3").Value = 0.0213
$ws.Range("G3").Value = 0.1036066459543428
$ws.Range("H3").Value = 0.1036066459543428
$ws.Range("I3").Value = 0.09253005538295286
$ws.Range("J3").Value = 0.07704742424509382
$ws.Range("K3").Value = 51.2
$ws.Range("L3").Value = 0.06916115088477645
$ws.Range("M3").Value = 23
$ws.Range("N3").Value = 0.0432574760203122
$ws.Range("O3").Value = 0.44921875
$ws.Range("P3").Value = 23
$ws.Range("Q3").Value = 0.0432574760203122
$ws.Range("R3").Value = 0.44921875
$ws.Range("U3").Value = 341
$ws.Range("V3").Value = 0.6413391009968027
$ws.Range("W3").Value = 0.1536614645858344
$ws.Range("X3").Value = 0.05334471444136124
$ws.Range("Y3").Value = 0.1003167501444731
$ws.Range("Z3").Value = 3.170449678800856
$ws.Range("AA3").Value = 0.244274981450291
$ws.Range("AB3").Value = 0.04847585350092396
$ws.Range("AC3").Value = 0.195799127949367
$ws.Range("AD3").Value = 146.2
$ws.Range("AF3").Value = 146.2
$ws.Range("AG3").Value = -194.8
$ws.Range("AH3").Value = 0.2156660274376751
$ws.Range("AI3").Value = 0.2491903869098346
$ws.Range("AJ3").Value = -0.5782131196200653
$ws.Range("AK3").Value = -0.7928367928367929
$ws.Range("AL3").Value = 9.9
$ws.Range("AM3").Value = 9.9
$ws.Range("AN3").Value = 1.973009446693657
$ws.Range("AO3").Value = 6.919191919191919
$ws.Range("AP3").Value = -2.628879892037787
$ws.Range("AQ3").Value = 6.919191919191919

# Row 4
$ws.Range("D4").Value = 0.00281
$ws.Range("E4").Value = 0.163
$ws.Range("G4").Value = 0.3566739606126915
$ws.Range("H4").Value = 0.3566739606126915
$ws.Range("I4").Value = 0.4048140043763676
$ws.Range("J4").Value = 0.4010313161387524
$ws.Range("K4").Value = 36.2
$ws.Range("L4").Value = 0.3960612691466083
$ws.Range("M4").Value = 16.4
$ws.Range("N4").Value = 0.07011543394613082
$ws.Range("O4").Value = 0.4530386740331491
$ws.Range("P4").Value = 16.4
$ws.Range("Q4").Value = 0.07011543394613082
$ws.Range("R4").Value = 0.4530386740331491
$ws.Range("U4").Value = 26.3
$ws.Range("V4").Value = 0.1124412141941
$ws.Range("W4").Value = 0.1114532019704434
$ws.Range("X4").Value = 0.04500067940374067
$ws.Range("Y4").Value = 0.06645252256670268
$ws.Range("Z4").Value = 0.3003614853762734
$ws.Range("AA4").Value = 0.1204543617978375
$ws.Range("AB4").Value = 0.04500067940374067
$ws.Range("AC4").Value = 0.07545368239409687
$ws.Range("AG4").Value = -26.3
$ws.Range("AJ4").Value = -0.1266859344894027
$ws.Range("AK4").Value = -0.08900169204737733
$ws.Range("AP4").Value = -0.7013333333333334

# Row 5
$ws.Range("D5").Value = 0.251
$ws.Range("E5").Value = 0.00749
$ws.Range("G5").Value = 0.1298654056020371
$ws.Range("H5").Value = 0.1298654056020371
$ws.Range("I5").Value = 0.1280465623863223
$ws.Range("J5").Value = 0.1267820229518106
$ws.Range("K5").Value = 37
$ws.Range("L5").Value = 0.1345943979628956
$ws.Range("M5").Value = 16.1
$ws.Range("N5").Value = 0.05048604578237693
$ws.Range("O5").Value = 0.4351351351351352
$ws.Range("P5").Value = 16.1
$ws.Range("Q5").Value = 0.05048604578237693
$ws.Range("R5").Value = 0.4351351351351352
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 26.3
$ws.Range("V5").Value = 0.08247099404201945
$ws.Range("W5").Value = 0.1072152999130687
$ws.Range("X5").Value = 0.04500067940374067
$ws.Range("Y5").Value = 0.06221462050932799
$ws.Range("Z5").Value = 0.8760356915232631
$ws.Range("AA5").Value = 0.1110655771493076
$ws.Range("AB5").Value = 0.04500067940374067
$ws.Range("AC5").Value = 0.06606489774556697
$ws.Range("AG5").Value = -26.3
$ws.Range("AJ5").Value = -0.08988380041011622
$ws.Range("AK5").Value = -0.07092772384034519
$ws.Range("AP5").Value = -0.7285318559556787

# Row 6
$ws.Range("D6").Value = -0.0155
$ws.Range("E6").Value = -0.0224
$ws.Range("G6").Value = 0.07009646302250803
$ws.Range("H6").Value = 0.07009646302250803
$ws.Range("I6").Value = 0.09694533762057878
$ws.Range("J6").Value = 0.09615008289790997
$ws.Range("K6").Value = 5.05
$ws.Range("L6").Value = 0.08118971061093247
$ws.Range("M6").Value = 3.53
$ws.Range("N6").Value = 0.07690631808278867
$ws.Range("O6").Value = 0.699009900990099
$ws.Range("P6").Value = 3.53
$ws.Range("Q6").Value = 0.07690631808278867
$ws.Range("R6").Value = 0.699009900990099
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 10.4
$ws.Range("V6").Value = 0.2265795206971678
$ws.Range("W6").Value = 0.04066022544283414
$ws.Range("X6").Value = 0.04500067940374067
$ws.Range("Y6").Value = -0.004340453960906529
$ws.Range("Z6").Value = 0.5159256801592568
$ws.Range("AA6").Value = 0.04960629691647313
$ws.Range("AB6").Value = 0.04500067940374067
$ws.Range("AC6").Value = 0.00460561751273246
$ws.Range("AG6").Value = -10.4
$ws.Range("AJ6").Value = -0.2929577464788732
$ws.Range("AK6").Value = -0.09885931558935362
$ws.Range("AP6").Value = -1.536189069423929

# Row 7
$ws.Range("D7").Value = 0.43
$ws.Range("I7").Value = 0.01175925925925926
$ws.Range("J7").Value = 0.01175925925925926
$ws.Range("K7").Value = -4.62
$ws.Range("L7").Value = -0.4277777777777778
$ws.Range("O7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("U7").Value = 0.049
$ws.Range("V7").Value = 0.006657608695652174
$ws.Range("W7").Value = -0.1196891191709845
$ws.Range("X7").Value = 0.04500067940374067
$ws.Range("Y7").Value = -0.1646897985747251
$ws.Range("Z7").Value = 0.2800321518396557
$ws.Range("AA7").Value = 0.003292970674410766
$ws.Range("AB7").Value = 0.04500067940374067
$ws.Range("AC7").Value = -0.04170770872932991
$ws.Range("AG7").Value = -0.049
$ws.Range("AJ7").Value = -0.006702229517165915
$ws.Range("AK7").Value = -0.001271043552696428

# Cells removed entirely (clear contents)
$ws.Range("E7").ClearContents()
